{"js": "// Update the worksheet date and all division problems to the new values.\n// Each old value in the document is unique, so a simple search/replace\n// per pair is safe and order-independent.\nconst replacements = [\n  [\"2024-07-19 Friday\", \"2024-07-20 Saturday\"],\n  [\"750\u00f73=\", \"168\u00f77=\"],\n  [\"815\u00f76=\", \"195\u00f78=\"],\n  [\"952\u00f73=\", \"140\u00f73=\"],\n  [\"611\u00f76=\", \"658\u00f76=\"],\n  [\"184\u00f72=\", \"849\u00f75=\"],\n  [\"369\u00f74=\", \"803\u00f79=\"],\n  [\"814\u00f75=\", \"616\u00f72=\"],\n  [\"542\u00f77=\", \"961\u00f75=\"],\n  [\"745\u00f75=\", \"554\u00f79=\"],\n  [\"260\u00f74=\", \"330\u00f75=\"],\n  [\"771\u00f74=\", \"946\u00f77=\"],\n  [\"682\u00f72=\", \"592\u00f76=\"],\n  [\"140\u00f75=\", \"156\u00f77=\"],\n  [\"553\u00f77=\", \"338\u00f76=\"],\n  [\"322\u00f77=\", \"894\u00f79=\"],\n  [\"111\u00f79=\", \"142\u00f78=\"],\n  [\"715\u00f73=\", \"395\u00f79=\"],\n  [\"388\u00f72=\", \"552\u00f79=\"],\n  [\"267\u00f79=\", \"832\u00f75=\"],\n  [\"379\u00f76=\", \"127\u00f76=\"],\n  [\"861\u00f72=\", \"644\u00f74=\"],\n  [\"991\u00f76=\", \"788\u00f75=\"],\n  [\"991\u00f79=\", \"887\u00f72=\"],\n  [\"186\u00f72=\", \"648\u00f73=\"],\n  [\"387\u00f73=\", \"675\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all division problems to the new values.\n# Each old value in the document is unique, so a simple Find/Replace per\n# pair (scoped to the whole document content) is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-19 Friday\", \"2024-07-20 Saturday\"),\n    @(\"750\u00f73=\", \"168\u00f77=\"),\n    @(\"815\u00f76=\", \"195\u00f78=\"),\n    @(\"952\u00f73=\", \"140\u00f73=\"),\n    @(\"611\u00f76=\", \"658\u00f76=\"),\n    @(\"184\u00f72=\", \"849\u00f75=\"),\n    @(\"369\u00f74=\", \"803\u00f79=\"),\n    @(\"814\u00f75=\", \"616\u00f72=\"),\n    @(\"542\u00f77=\", \"961\u00f75=\"),\n    @(\"745\u00f75=\", \"554\u00f79=\"),\n    @(\"260\u00f74=\", \"330\u00f75=\"),\n    @(\"771\u00f74=\", \"946\u00f77=\"),\n    @(\"682\u00f72=\", \"592\u00f76=\"),\n    @(\"140\u00f75=\", \"156\u00f77=\"),\n    @(\"553\u00f77=\", \"338\u00f76=\"),\n    @(\"322\u00f77=\", \"894\u00f79=\"),\n    @(\"111\u00f79=\", \"142\u00f78=\"),\n    @(\"715\u00f73=\", \"395\u00f79=\"),\n    @(\"388\u00f72=\", \"552\u00f79=\"),\n    @(\"267\u00f79=\", \"832\u00f75=\"),\n    @(\"379\u00f76=\", \"127\u00f76=\"),\n    @(\"861\u00f72=\", \"644\u00f74=\"),\n    @(\"991\u00f76=\", \"788\u00f75=\"),\n    @(\"991\u00f79=\", \"887\u00f72=\"),\n    @(\"186\u00f72=\", \"648\u00f73=\"),\n    @(\"387\u00f73=\", \"675\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
